$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers: insert an extra "N" in "N2O" -> "NN2O"
$ws.Range("C1").Value = "avg_Gasera_NN2O_flux_mgm2h"
$ws.Range("D1").Value = "avg_Gasera_NN2O_flux_mgm2h_cor"

# Apply new transformed values for avg_Gasera_NN2O_flux_mgm2h (column C, rows 2-18)
$newC = @{
    2  = -0.04550788477675025
    3  = -0.2012796706260847
    4  = -0.0006252952788222648
    5  = 0.07720773154899609
    6  = 0.1382331609005755
    7  = 0.04429885970070662
    8  = 0.009137731916897571
    9  = -0.05625733464332385
    10 = 0.004770541528519517
    11 = 0.01974554801672977
    12 = 0.01631969987381054
    13 = -0.006149331825670259
    14 = 0.04682481208409594
    15 = -0.01976378274983011
    16 = 0.005420904679603111
    17 = 0.01646375253003778
    18 = 0.01011357495795964
}

foreach ($row in $newC.Keys) {
    $ws.Cells.Item($row, 3).Value = $newC[$row]
}

# Update the corrected value for row 5 (avg_Gasera_NN2O_flux_mgm2h_cor, column D)
$ws.Range("D5").Value = 0.08000959782295058
